$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table5: show the totals row (extends ref from G9:J10 to G9:J11) ---
$table5 = $ws.ListObjects.Item("Table5")
$table5.ShowTotals = $true

# --- Row 24 / Row 25 (new "Table7" source data) ---
# Write header text and values in the same left-to-right / top-to-bottom order the
# original author appears to have used, skipping the D column header for now so the
# shared-string table indices land in the right slots.
$ws.Range("B24").Value = "1 Rotation distance (mm)"
$ws.Range("C24").Value = "Inner-Radius (mm)"
$ws.Range("E24").Value = "Rotations per mL"
$ws.Range("F24").Value = "Volume per Day (mL)"
$ws.Range("G24").Value = "Rotations requered per day"

$ws.Range("B25").Value = 1.25
$ws.Range("C25").Value = 7.95
$ws.Range("F25").Value = 2.5

# --- Row 27 / Row 28 (new "Table8" source data) + the late-added D24 header ---
$ws.Range("D27").Value = "Volume per Rotation (mL/rot)"
$ws.Range("D24").Value = "mm per mL"
$ws.Range("E27").Value = "Rotations per Day"
$ws.Range("F27").Value = "motor rpm"
$ws.Range("G27").Value = "steps per sec"
$ws.Range("H27").Value = "microsteps period"
$ws.Range("I27").Value = "Constant"

# --- Formulas for row 25 (Table7 body) ---
$ws.Range("D25").Formula = "=1000/(3.14159*(POWER(C25,2)))"
$ws.Range("E25").Formula = "=D25/B25"
$ws.Range("G25").Formula = "=E25*F25"

# --- Formulas for row 28 (Table8 body) ---
$ws.Range("D28").Formula = "=B25/D25"
$ws.Range("E28").Formula = "=F25/D28"
$ws.Range("F28").Formula = "=E28/1440*20"
$ws.Range("G28").Formula = "=(F28*K3)/60"
$ws.Range("H28").Formula = "=G28/C14"
$ws.Range("I28").Formula = "=H28/F25"

# --- Turn the two new ranges into tables (Excel auto-assigns Table7 / Table8) ---
$table7 = $ws.ListObjects.Add(1, $ws.Range("B24:G25"), 0, 1)
$table7.TableStyle = "TableStyleLight2"

$table8 = $ws.ListObjects.Add(1, $ws.Range("D27:I28"), 0, 1)
$table8.TableStyle = "TableStyleLight2"

# --- Cosmetic touch-ups called out by the diff ---
$ws.Columns.Item(2).ColumnWidth = 22.0
$ws.Columns.Item(4).ColumnWidth = 25.6
$ws.Columns.Item(6).ColumnWidth = 18.7

$ws.Range("C26").Select()
